# Update threshold values on the FINE and COARSE sheets (row 3)
$wb = $excel.ActiveWorkbook

# FINE sheet: row 3, columns C-F
$wsFine = $wb.Worksheets.Item("FINE")
$wsFine.Range("C3").Value = 0.25
$wsFine.Range("D3").Value = 0.37
$wsFine.Range("E3").Value = 0.47
$wsFine.Range("F3").Value = 0.57

# COARSE sheet: row 3, column D
$wsCoarse = $wb.Worksheets.Item("COARSE")
$wsCoarse.Range("D3").Value = 0.5
